# Car sheet ("汽車") gains explicit column headers plus the extra
# property_category / category / date / legislator_name / legislator_id /
# source_file / index columns that the other property sheets already carry,
# along with a new "capacity" header for the (until now unlabeled)
# engine-capacity value that already lived in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # 汽車 (car) sheet

# --- Row 1: proper header labels -------------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2: fill in the newly added data columns ---------------------------------
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# Leading apostrophe keeps this a text value instead of letting Excel coerce
# it into a date serial; the formatting pass below normalizes the style.
$ws.Range("J2").Value = "'2012-04-25"
$ws.Range("K2").Value = "柯建銘"
$ws.Range("L2").Value = 629
$ws.Range("M2").Value = "tmp81dc1"
$ws.Range("N2").Value = 29

# --- Match formatting of the existing header/data cells --------------------------
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0
